$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.603.00"
$ws.Range("E2").Value = "  +4.13%  "
$ws.Range("D3").Value = "1.744.82"
$ws.Range("E3").Value = "  +4.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.23"
$ws.Range("E5").Value = "  +3.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4814"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2699"
$ws.Range("E8").Value = "  +3.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06266"
$ws.Range("E9").Value = "  +1.61%  "
$ws.Range("D10").Value = "1.744.17"
$ws.Range("E10").Value = "  +4.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07121"
$ws.Range("E11").Value = "  +1.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.86"
$ws.Range("E12").Value = "  +7.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6238"
$ws.Range("E13").Value = "  +6.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.515"
$ws.Range("E14").Value = "  +3.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.56"
$ws.Range("E15").Value = "  +3.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "26.586.09"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006912"
$ws.Range("E19").Value = "  +2.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.73"
$ws.Range("E20").Value = "  +2.69%  "
$ws.Range("D21").Value = "1.966.90"
$ws.Range("E21").Value = "  +4.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.637"
$ws.Range("E22").Value = "  +4.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.835"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.362"
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.95"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.41"
$ws.Range("E26").Value = "  +2.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.822"
$ws.Range("E27").Value = "  +5.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.436"
$ws.Range("E28").Value = "  +3.56%  "
$ws.Range("E29").Value = "  +2.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.013"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.758"
$ws.Range("E31").Value = "  +3.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07887"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04612"
$ws.Range("E33").Value = "  +7.89%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.619"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6459"
$ws.Range("E35").Value = "  +6.60%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9994"
$ws.Range("E36").Value = "  +4.68%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9484"
$ws.Range("E37").Value = "  +6.14%  "
$ws.Range("B38").Value = "Quant"
$ws.Range("C38").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "113.52"
$ws.Range("E38").Value = "  +18.33%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.003"
$ws.Range("E39").Value = "  +7.94%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.433"
$ws.Range("E40").Value = "  -6.13%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.003"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.770"
$ws.Range("E42").Value = "  +18.17%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.01512"
$ws.Range("E43").Value = "  +2.37%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3927"
$ws.Range("E44").Value = "  +4.71%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1214"
$ws.Range("E45").Value = "  +8.72%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.732"
$ws.Range("E46").Value = "  +8.51%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05329"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.998"
$ws.Range("E48").Value = "  +7.67%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.77"
$ws.Range("E49").Value = "  +2.73%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.272"
$ws.Range("E50").Value = "  +5.88%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3456"
$ws.Range("E51").Value = "  +3.63%  "
